# Generate Report for handback
#
# Marks the two pending-handoff files (94a3b0c1... and 9be84bbd...) as
# "Handed back" on both the zh-cn and de-de localization-status sheets,
# fills in the now-known "Latest Target File" / "Latest Handback File"
# columns (E/F) with hyperlinks matching the originally handed-off file,
# and stamps the "Latest Handback DateTime" column (G) with the actual
# handback timestamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Status column: no longer pending handoff, it has now been handed back.
$ws.Range("B2").Value = "Handed back"
$ws.Range("B3").Value = "Handed back"

# Row 2 (94a3b0c1-3797-4dd1-b5f2-dc182c966830.md)
$ws.Range("E2").Value = "94a3b0c1-3797-4dd1-b5f2-dc182c966830.md"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/e8842656f23126b1e2bb6385f0dc0dd2198b2334/e2e/94a3b0c1-3797-4dd1-b5f2-dc182c966830.md", "", "", "94a3b0c1-3797-4dd1-b5f2-dc182c966830.md") | Out-Null

$ws.Range("F2").Value = "94a3b0c1-3797-4dd1-b5f2-dc182c966830.47062af9c8ed34cc7f3af68ab32df37d279a283f.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/796703b5262d5067ecbd06a53701d351fb7bbdcf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/94a3b0c1-3797-4dd1-b5f2-dc182c966830.47062af9c8ed34cc7f3af68ab32df37d279a283f.zh-cn.xlf", "", "", "94a3b0c1-3797-4dd1-b5f2-dc182c966830.47062af9c8ed34cc7f3af68ab32df37d279a283f.zh-cn.xlf") | Out-Null

$ws.Range("G2").Value = "2016-01-07 12:55:26"

# Row 3 (9be84bbd-15b2-494b-b15b-281c65f3c082.md)
$ws.Range("E3").Value = "9be84bbd-15b2-494b-b15b-281c65f3c082.md"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/e8842656f23126b1e2bb6385f0dc0dd2198b2334/e2e/9be84bbd-15b2-494b-b15b-281c65f3c082.md", "", "", "9be84bbd-15b2-494b-b15b-281c65f3c082.md") | Out-Null

$ws.Range("F3").Value = "9be84bbd-15b2-494b-b15b-281c65f3c082.df876eee37292aaa9718efe6365569a9768f969c.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/796703b5262d5067ecbd06a53701d351fb7bbdcf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/9be84bbd-15b2-494b-b15b-281c65f3c082.df876eee37292aaa9718efe6365569a9768f969c.zh-cn.xlf", "", "", "9be84bbd-15b2-494b-b15b-281c65f3c082.df876eee37292aaa9718efe6365569a9768f969c.zh-cn.xlf") | Out-Null

$ws.Range("G3").Value = "2016-01-07 12:55:26"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = "Handed back"
$ws.Range("B3").Value = "Handed back"

# Row 2 (94a3b0c1-3797-4dd1-b5f2-dc182c966830.md)
$ws.Range("E2").Value = "94a3b0c1-3797-4dd1-b5f2-dc182c966830.md"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/e8842656f23126b1e2bb6385f0dc0dd2198b2334/e2e/94a3b0c1-3797-4dd1-b5f2-dc182c966830.md", "", "", "94a3b0c1-3797-4dd1-b5f2-dc182c966830.md") | Out-Null

$ws.Range("F2").Value = "94a3b0c1-3797-4dd1-b5f2-dc182c966830.47062af9c8ed34cc7f3af68ab32df37d279a283f.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9a5515c19deafc835c895fc6654a17f238351fe0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/94a3b0c1-3797-4dd1-b5f2-dc182c966830.47062af9c8ed34cc7f3af68ab32df37d279a283f.de-de.xlf", "", "", "94a3b0c1-3797-4dd1-b5f2-dc182c966830.47062af9c8ed34cc7f3af68ab32df37d279a283f.de-de.xlf") | Out-Null

$ws.Range("G2").Value = "2016-01-07 12:55:44"

# Row 3 (9be84bbd-15b2-494b-b15b-281c65f3c082.md)
$ws.Range("E3").Value = "9be84bbd-15b2-494b-b15b-281c65f3c082.md"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/e8842656f23126b1e2bb6385f0dc0dd2198b2334/e2e/9be84bbd-15b2-494b-b15b-281c65f3c082.md", "", "", "9be84bbd-15b2-494b-b15b-281c65f3c082.md") | Out-Null

$ws.Range("F3").Value = "9be84bbd-15b2-494b-b15b-281c65f3c082.df876eee37292aaa9718efe6365569a9768f969c.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9a5515c19deafc835c895fc6654a17f238351fe0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/9be84bbd-15b2-494b-b15b-281c65f3c082.df876eee37292aaa9718efe6365569a9768f969c.de-de.xlf", "", "", "9be84bbd-15b2-494b-b15b-281c65f3c082.df876eee37292aaa9718efe6365569a9768f969c.de-de.xlf") | Out-Null

$ws.Range("G3").Value = "2016-01-07 12:55:44"
